$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33470
$ws.Range("B3").Value = 71711
$ws.Range("B4").Value = 47137
$ws.Range("B5").Value = 712
$ws.Range("B6").Value = 1369
$ws.Range("B7").Value = 127
$ws.Range("B8").Value = 782
